# Backup QR Scanner data - append the latest scan log entry as a new row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores every value as text (Student ID, dates and times are all
# kept as plain strings), so force the new row to a text number format
# before assigning the values. This stops Excel from auto-converting the
# numeric-looking strings ("235166") or the date/time strings into real
# numbers/dates.
$newRow = 3
$ws.Range("A$newRow`:F$newRow").NumberFormat = "@"

$ws.Range("A$newRow").Value = "235166"
$ws.Range("B$newRow").Value = "Parasitology SGD/POS"
$ws.Range("C$newRow").Value = "12/10/2025"
$ws.Range("D$newRow").Value = "15:13:46"
$ws.Range("E$newRow").Value = "Manual"
$ws.Range("F$newRow").Value = "Mayarembaby@med.asu.edu.eg"

# Mark the "number stored as text" warning as ignored for the newly used
# range, matching the expanded ignoredErrors sqref (A1:F3) for the sheet.
$ws.Range("A1:F$newRow").Errors().Item(9).Ignore = $true
